# WB_Test_Report_2025-12-26.xlsx update
# - Rename Package_1 -> BOP_1, Package_2 -> BOP_2
# - Add three new milestone-detail sheets: BOP_3, BOP_4, BOP_5
# - Refresh the Summary sheet with five BOP iterations (was two Package iterations)
# - Refresh milestone timelines on BOP_1 / BOP_2 and populate BOP_3..BOP_5

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Helper: write a value, forcing "text" storage (quote-prefix) when the
# literal looks numeric so it round-trips the same way the source data does.
function Set-TextValue {
    param($cell, [string]$text)
    if ($text -match '^-?[0-9]+(\.[0-9]+)?$') {
        $cell.Value = "'" + $text
    } else {
        $cell.Value = $text
    }
}

# ---------------------------------------------------------------------------
# 1) Rename the two existing package sheets to the new BOP naming
# ---------------------------------------------------------------------------
$wsBop1 = $wb.Worksheets.Item("Package_1")
$wsBop1.Name = "BOP_1"

$wsBop2 = $wb.Worksheets.Item("Package_2")
$wsBop2.Name = "BOP_2"

# ---------------------------------------------------------------------------
# 2) Add three new sheets (BOP_3, BOP_4, BOP_5) at the end of the workbook
# ---------------------------------------------------------------------------
$newNames = @("BOP_3", "BOP_4", "BOP_5")
foreach ($name in $newNames) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet = $wb.Worksheets.Add($null, $lastSheet)
    $newSheet.Name = $name
}

# ---------------------------------------------------------------------------
# 3) Summary sheet: five BOP iterations
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$summaryRows = @(
    @(1, "BOP (MI)", "N/A",        "N/A",        "FAILED", "134.88", "MI"),
    @(2, "BOP (WI)", "3003179758", "1003052794", "PASSED", "464.00", "WI"),
    @(3, "BOP (DE)", "3003179757", "1003052797", "PASSED", "473.12", "DE"),
    @(4, "BOP (PA)", "3003179759", "1003052795", "PASSED", "474.18", "PA"),
    @(5, "BOP (OH)", "3003179760", "1003052796", "PASSED", "468.78", "OH")
)

$r = 2
foreach ($row in $summaryRows) {
    $wsSummary.Cells.Item($r, 1).Value = $row[0]
    Set-TextValue $wsSummary.Cells.Item($r, 2) $row[1]
    Set-TextValue $wsSummary.Cells.Item($r, 3) $row[2]
    Set-TextValue $wsSummary.Cells.Item($r, 4) $row[3]
    Set-TextValue $wsSummary.Cells.Item($r, 5) $row[4]
    Set-TextValue $wsSummary.Cells.Item($r, 6) $row[5]
    Set-TextValue $wsSummary.Cells.Item($r, 7) $row[6]
    $r++
}

# ---------------------------------------------------------------------------
# 4) BOP_1 (was Package_1): only the "Account Created" / "Test Execution
#    Failed" duration + timestamp cells moved
# ---------------------------------------------------------------------------
$wsBop1.Cells.Item(2, 3).Value = "69.30s"
$wsBop1.Cells.Item(2, 4).Value = "2025-12-26T14:54:21.518Z"
$wsBop1.Cells.Item(3, 3).Value = "65.58s"
$wsBop1.Cells.Item(3, 4).Value = "2025-12-26T14:55:27.106Z"

# ---------------------------------------------------------------------------
# 5) BOP_2 (was Package_2) + BOP_3 / BOP_4 / BOP_5: a full seven-row
#    milestone timeline (header + 6 milestones)
# ---------------------------------------------------------------------------
$header = @("Milestone", "Status", "Duration (s)", "Timestamp")

$milestoneData = @{
    "BOP_2" = @(
        @("Account Created",                     "PASSED", "69.29s",  "2025-12-26T14:54:16.289Z"),
        @("Building and Classification Added",    "PASSED", "142.62s", "2025-12-26T14:56:38.914Z"),
        @("Quote Rated Successfully",              "PASSED", "22.45s",  "2025-12-26T14:57:01.369Z"),
        @("Submitting for Approval",                "PASSED", "33.98s",  "2025-12-26T14:57:35.356Z"),
        @("UW Issues Approved in PolicyCenter",      "PASSED", "44.90s",  "2025-12-26T14:58:20.253Z"),
        @("Policy Issued Successfully",               "PASSED", "150.76s", "2025-12-26T15:00:51.014Z")
    )
    "BOP_3" = @(
        @("Account Created",                     "PASSED", "71.25s",  "2025-12-26T14:54:15.823Z"),
        @("Building and Classification Added",    "PASSED", "150.96s", "2025-12-26T14:56:46.787Z"),
        @("Quote Rated Successfully",              "PASSED", "23.44s",  "2025-12-26T14:57:10.233Z"),
        @("Submitting for Approval",                "PASSED", "35.31s",  "2025-12-26T14:57:45.551Z"),
        @("UW Issues Approved in PolicyCenter",      "PASSED", "41.13s",  "2025-12-26T14:58:26.685Z"),
        @("Policy Issued Successfully",               "PASSED", "151.03s", "2025-12-26T15:00:57.717Z")
    )
    "BOP_4" = @(
        @("Account Created",                     "PASSED", "71.17s",  "2025-12-26T14:54:16.685Z"),
        @("Building and Classification Added",    "PASSED", "150.67s", "2025-12-26T14:56:47.355Z"),
        @("Quote Rated Successfully",              "PASSED", "24.13s",  "2025-12-26T14:57:11.486Z"),
        @("Submitting for Approval",                "PASSED", "34.10s",  "2025-12-26T14:57:45.584Z"),
        @("UW Issues Approved in PolicyCenter",      "PASSED", "41.95s",  "2025-12-26T14:58:27.539Z"),
        @("Policy Issued Successfully",               "PASSED", "152.16s", "2025-12-26T15:00:59.697Z")
    )
    "BOP_5" = @(
        @("Account Created",                     "PASSED", "69.89s",  "2025-12-26T14:54:21.586Z"),
        @("Building and Classification Added",    "PASSED", "140.51s", "2025-12-26T14:56:42.103Z"),
        @("Quote Rated Successfully",              "PASSED", "32.18s",  "2025-12-26T14:57:14.288Z"),
        @("Submitting for Approval",                "PASSED", "33.88s",  "2025-12-26T14:57:48.170Z"),
        @("UW Issues Approved in PolicyCenter",      "PASSED", "40.41s",  "2025-12-26T14:58:28.587Z"),
        @("Policy Issued Successfully",               "PASSED", "151.91s", "2025-12-26T15:01:00.497Z")
    )
}

foreach ($sheetName in @("BOP_2", "BOP_3", "BOP_4", "BOP_5")) {
    $ws = $wb.Worksheets.Item($sheetName)

    for ($c = 1; $c -le 4; $c++) {
        $ws.Cells.Item(1, $c).Value = $header[$c - 1]
    }

    $rowIdx = 2
    foreach ($milestone in $milestoneData[$sheetName]) {
        for ($c = 1; $c -le 4; $c++) {
            $ws.Cells.Item($rowIdx, $c).Value = $milestone[$c - 1]
        }
        $rowIdx++
    }
}

Write-Host "Workbook updated: Summary + BOP_1..BOP_5"
